# Add the "Final Demo" worksheet as the last tab, matching the structure
# used by the existing "Section 1" / "Section 2" question sheets.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Final Demo"

# Match the column sizing used on the other question sheets (A:P width 29).
$newSheet.Columns("A:P").ColumnWidth = 28.14

# ---- Question 4 (col G/H): "How Long Does a Typical New Year's Resolution
# Last?" -- answers are keyed in collection order, then the cells get
# resorted by popularity (hence the scratch-then-final-fill below). ----
$newSheet.Range("G1").Value = "How Long Does a Typical New Year's Resolution Last?"

$newSheet.Range("Z1").Value = "1 Month"
$newSheet.Range("Z2").Value = "1 Week"
$newSheet.Range("Z3").Value = "3 Months"
$newSheet.Range("Z4").Value = "1 Year"
$newSheet.Range("Z5").Value = "2 Weeks"
$newSheet.Range("Z1:Z5").ClearContents()

$newSheet.Range("G2").Value = "1 Month"
$newSheet.Range("H2").Value = 41
$newSheet.Range("G3").Value = "1 Week"
$newSheet.Range("H3").Value = 30
$newSheet.Range("G4").Value = "2 Weeks"
$newSheet.Range("H4").Value = 10
$newSheet.Range("G5").Value = "1 Year"
$newSheet.Range("H5").Value = 5
$newSheet.Range("G6").Value = "3 Months"
$newSheet.Range("H6").Value = 3

# ---- Question 1 (col A/B): "Name a Popular Chirstmas Carol" ----
$newSheet.Range("A1").Value = "Name a Popular Chirstmas Carol"
$newSheet.Range("A2").Value = '"We Wish You a Merry Christmas"'
$newSheet.Range("B2").Value = 21
$newSheet.Range("A3").Value = '"Jingle Bells"'
$newSheet.Range("B3").Value = 18
$newSheet.Range("A4").Value = '"Frosty the Snowman"'
$newSheet.Range("B4").Value = 15
$newSheet.Range("A5").Value = '"Deck the Halls"'
$newSheet.Range("B5").Value = 10
$newSheet.Range("A6").Value = '"Joy to the World"'
$newSheet.Range("B6").Value = 8
$newSheet.Range("A7").Value = "`"Santa Claus Is Comin' to Town`""
$newSheet.Range("B7").Value = 7
$newSheet.Range("A8").Value = '"Holly Jolly Christmas"'
$newSheet.Range("B8").Value = 7
$newSheet.Range("A9").Value = '"Silent Night"'
$newSheet.Range("B9").Value = 5

# ---- Question 2 (col C/D): "Name an Item You'd Need to Dress Up as Santa
# Clause" ----
$newSheet.Range("C1").Value = "Name an Item You'd Need to Dress Up as Santa Clause"
$newSheet.Range("C2").Value = "Santa Hat"
$newSheet.Range("D2").Value = 40
$newSheet.Range("C3").Value = "Beard"
$newSheet.Range("D3").Value = 24
$newSheet.Range("C4").Value = "Big Belly"
$newSheet.Range("D4").Value = 9
$newSheet.Range("C5").Value = "Red Coat"
$newSheet.Range("D5").Value = 8
$newSheet.Range("C6").Value = "Red Pants/Belt"
$newSheet.Range("D6").Value = 5
$newSheet.Range("C7").Value = "Suspenders"
$newSheet.Range("D7").Value = 2
$newSheet.Range("C8").Value = "Black Boots"
$newSheet.Range("D8").Value = 2

# ---- Question 3 (col E/F): "Name One of Santa's Reindeer" ----
$newSheet.Range("E1").Value = "Name One of Santa's Reindeer"
$newSheet.Range("E2").Value = "Rudolph"
$newSheet.Range("F2").Value = 32
$newSheet.Range("E3").Value = "Comet"
$newSheet.Range("F3").Value = 14
$newSheet.Range("E4").Value = "Donnor"
$newSheet.Range("F4").Value = 12
$newSheet.Range("E5").Value = "Blitzen"
$newSheet.Range("F5").Value = 11
$newSheet.Range("E6").Value = "Cupid"
$newSheet.Range("F6").Value = 9
$newSheet.Range("E7").Value = "Prancer"
$newSheet.Range("F7").Value = 6
$newSheet.Range("E8").Value = "Dasher"
$newSheet.Range("F8").Value = 4
$newSheet.Range("E9").Value = "Vixen"
$newSheet.Range("F9").Value = 3

# Center-align the question header row, like the other sheets.
$newSheet.Range("A1:H1").HorizontalAlignment = -4108

# Apply the data-row font formatting (size 11, black) the other sheets use,
# scoped to exactly the filled rows of each answer column (so we don't
# stamp empty styled cells beyond the real data, like Excel wouldn't).
foreach ($rng in @("A2:B9", "C2:D8", "E2:F9", "G2:H6")) {
    $r = $newSheet.Range($rng)
    $r.Font.Size = 11
    $r.Font.Color = 0
}

# Left-over formatted-but-empty cells carried over from the template (as in
# the other question sheets' trailing, now-unused answer columns).
foreach ($addr in @("G7", "H7", "J7", "L7", "G8", "H8", "J8", "L8", "L9", "L10")) {
    $c = $newSheet.Range($addr)
    $c.Font.Size = 11
    $c.Font.Color = 0
}
$newSheet.Range("O7").NumberFormat = "0%"

# Header/data cell merges, matching the other question sheets.
$newSheet.Range("A1:B1").Merge()
$newSheet.Range("C1:D1").Merge()
$newSheet.Range("E1:F1").Merge()
$newSheet.Range("G1:H1").Merge()

$newSheet.Range("A1:B1").Select()

# Make the new sheet the active/visible tab, like it was when last saved.
$newSheet.Activate()
